$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns with latest values.
# Force text format on Price cells before writing so strings such as "5.20" or
# "37.314.00" are preserved exactly (not auto-converted to numbers by Excel).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.314.00'
$ws.Range("E2").Value = '  +2.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.091.67'
$ws.Range("E3").Value = '  +3.95%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.12'
$ws.Range("E5").Value = '  +2.69%  '

$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.38'
$ws.Range("E8").Value = '  +22.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.82'
$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("E10").Value = '  +3.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0743'
$ws.Range("E11").Value = '  +4.04%  '

$ws.Range("E12").Value = '  +8.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.09'
$ws.Range("E13").Value = '  +5.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.397.30'
$ws.Range("E14").Value = '  +4.17%  '

$ws.Range("E15").Value = '  +4.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.090.37'
$ws.Range("E16").Value = '  +3.73%  '

$ws.Range("E17").Value = '  +6.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.247.56'
$ws.Range("E18").Value = '  +2.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.78'
$ws.Range("E19").Value = '  +2.34%  '

$ws.Range("E20").Value = '  +14.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0847'
$ws.Range("E21").Value = '  +4.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.92'
$ws.Range("E22").Value = '  +1.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.20'
$ws.Range("E23").Value = '  +6.61%  '

$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +2.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.33'
$ws.Range("E26").Value = '  +4.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.22'
$ws.Range("E27").Value = '  +7.05%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.70'
$ws.Range("E28").Value = '  +6.02%  '

$ws.Range("E29").Value = '  +2.86%  '

$ws.Range("E30").Value = '  +1.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.66'
$ws.Range("E31").Value = '  +9.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.09'
$ws.Range("E32").Value = '  +28.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.52'
$ws.Range("E33").Value = '  +4.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0617'
$ws.Range("E34").Value = '  +6.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0921'
$ws.Range("E35").Value = '  +6.01%  '

$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("E37").Value = '  +4.21%  '

$ws.Range("E38").Value = '  -0.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.25'
$ws.Range("E39").Value = '  +5.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.35'
$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("E41").Value = '  +15.58%  '

$ws.Range("E42").Value = '  +5.93%  '

$ws.Range("E43").Value = '  +6.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.72'
$ws.Range("E44").Value = '  +3.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0936'
$ws.Range("E45").Value = '  +15.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.30'
$ws.Range("E46").Value = '  +118.73%  '

$ws.Range("E47").Value = '  +2.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.323.30'
$ws.Range("E48").Value = '  +1.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.94'
$ws.Range("E49").Value = '  +6.03%  '

$ws.Range("E50").Value = '  +7.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.01'
$ws.Range("E51").Value = '  +15.18%  '
